# Update scripts with new TPM values (Il15-Il2rb LR-pair table).
#
# The underlying NATMI computation was re-run, which changed several
# numeric columns and re-assigned the Target cluster (column D) /
# Sending cluster (column A) values for each row, and reduced the number
# of sending/target-cluster combinations from 6 to 3 (rows 5-7 removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ECs -> FAPs) ---------------------------------------------------
$ws.Range("D2").Value2 = "FAPs"
$ws.Range("G2").Value2 = 4.922976999999999
$ws.Range("H2").Value2 = 14.768931
$ws.Range("I2").Value2 = 0.2897120038548413
$ws.Range("J2").Value2 = 0.2897120038548412
$ws.Range("K2").Value2 = 2
$ws.Range("L2").Value2 = 0.6666666666666666
$ws.Range("M2").Value2 = 0.092254
$ws.Range("N2").Value2 = 0.276762
$ws.Range("O2").Value2 = 1
$ws.Range("P2").Value2 = 1
$ws.Range("Q2").Value2 = 0.454164320158
$ws.Range("R2").Value2 = 4.087478881421999
$ws.Range("S2").Value2 = 0.2897120038548413
$ws.Range("T2").Value2 = 0.2897120038548412

# --- Row 3 (ECs -> FAPs sending cluster) -----------------------------------
$ws.Range("A3").Value2 = "FAPs"
$ws.Range("G3").Value2 = 9.803896
$ws.Range("H3").Value2 = 29.411688
$ws.Range("I3").Value2 = 0.5769489387710858
$ws.Range("J3").Value2 = 0.5769489387710858
$ws.Range("O3").Value2 = 1
$ws.Range("P3").Value2 = 1
$ws.Range("Q3").Value2 = 0.904448621584
$ws.Range("R3").Value2 = 8.140037594256
$ws.Range("S3").Value2 = 0.5769489387710858
$ws.Range("T3").Value2 = 0.5769489387710858

# --- Row 4 (FAPs -> MuSCs sending cluster, ECs -> FAPs target) ------------
$ws.Range("A4").Value2 = "MuSCs"
$ws.Range("D4").Value2 = "FAPs"
$ws.Range("G4").Value2 = 2.265785
$ws.Range("H4").Value2 = 6.797355
$ws.Range("I4").Value2 = 0.133339057374073
$ws.Range("J4").Value2 = 0.133339057374073
$ws.Range("K4").Value2 = 2
$ws.Range("L4").Value2 = 0.6666666666666666
$ws.Range("M4").Value2 = 0.092254
$ws.Range("N4").Value2 = 0.276762
$ws.Range("O4").Value2 = 1
$ws.Range("P4").Value2 = 1
$ws.Range("Q4").Value2 = 0.20902772939
$ws.Range("R4").Value2 = 1.88124956451
$ws.Range("S4").Value2 = 0.133339057374073
$ws.Range("T4").Value2 = 0.133339057374073

# --- Remove now-obsolete rows 5-7 and shift the remainder up --------------
$ws.Range("A5:T7").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp) | Out-Null
